# "add practical side projects"
#
# Under the "Side Projects" heading the paragraph that used to read
# "Exploration" becomes "Practical", and two new description paragraphs
# (Private Gitolite Server / Virtualized Dev Server) are inserted right
# after it, followed by a new "Tech Demos" heading paragraph that now
# introduces the pre-existing Imgboard/Arimaa entries.

$d = $word.ActiveDocument

# Helper: insert literal text at an absolute document position and
# return the position immediately following the inserted text so callers
# can chain inserts left-to-right without relying on a live Range object
# (collapsing a Range that sits exactly on a paragraph-mark boundary is
# ambiguous between "end of this paragraph" and "start of the next").
function InsertAt($pos, $text) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    return $pos + $text.Length
}

# Locate the "Exploration" paragraph under Side Projects.
$explorationIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd("`r") -eq "Exploration") {
        $explorationIdx = $i
        break
    }
}
if ($explorationIdx -eq 0) {
    throw "Could not find the 'Exploration' paragraph"
}

# Rename it to "Practical" (modifies the existing run in place).
$d.Paragraphs($explorationIdx).Range.Text = "Practical"
$practicalIdx = $explorationIdx

# --- Insert the "Private Gitolite Server" paragraph right after it ----
$rng = $d.Paragraphs($practicalIdx).Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$gitoliteIdx = $practicalIdx + 1

$pos = $d.Paragraphs($gitoliteIdx).Range.Start
$pos = InsertAt $pos ":"
$pos = InsertAt $pos " "
$gitoliteItalicStart = $pos
$pos = InsertAt $pos "Private Gitolite Server"
$gitoliteItalicEnd = $pos
$pos = InsertAt $pos " "
$pos = InsertAt $pos "- Ubuntu Server, gitolite, ssh"
$pos = InsertAt $pos " "
$pos = InsertAt $pos " "
$pos = InsertAt $pos "- Hosted on Rackspace"

# --- Insert the "Virtualized Dev Server" paragraph after that ---------
$rng2 = $d.Paragraphs($gitoliteIdx).Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
$devServerIdx = $gitoliteIdx + 1

$pos = $d.Paragraphs($devServerIdx).Range.Start
$pos = InsertAt $pos ":"
$pos = InsertAt $pos " "
$devServerItalicStart = $pos
$pos = InsertAt $pos "Virtualized Dev Server"
$devServerItalicEnd = $pos
$pos = InsertAt $pos " "
$pos = InsertAt $pos "- Ubuntu Server, KVM, libvert"
$pos = InsertAt $pos " "
$pos = InsertAt $pos " "
$pos = InsertAt $pos "- Hosted on local hardware"

# --- Insert the new "Tech Demos" sub-heading paragraph ----------------
$rng3 = $d.Paragraphs($devServerIdx).Range
$rng3.Collapse(0)
$rng3.InsertParagraphAfter()
$techDemosIdx = $devServerIdx + 1
$d.Paragraphs($techDemosIdx).Range.Text = "Tech Demos"

# Italicize the two project-name runs now that every paragraph break has
# already been created (doing this earlier makes freshly inserted
# paragraph marks/runs inherit the italic formatting).
$d.Range($gitoliteItalicStart, $gitoliteItalicEnd).Font.Italic = 1
$d.Range($devServerItalicStart, $devServerItalicEnd).Font.Italic = 1

Write-Host "Renamed paragraph $practicalIdx to 'Practical'; inserted 'Private Gitolite Server', 'Virtualized Dev Server' and 'Tech Demos' paragraphs."
